# PoC v1.2.0 - Versão final da Prova de Conceito para deployment
#
# This script reproduces the authored edit on Recursos.xlsx:
#   1. Duplicate the "machine" sheet into a new "machine (2)" sheet placed
#      right after it (preserves the original 3-machine data set).
#   2. On the original "machine" sheet, remove the "ASHE2 / Sliter 2" row
#      (row 3), which shifts "ATLAS / Sliter 3" up from row 4 to row 3.
#   3. Update the CalendarioId on the remaining two rows from
#      "CAL-PADRAO-5x8" to the new calendar "CAL-24x5".
#   4. Leave the first sheet's selection on the (now last) data row, and
#      re-activate it so it stays the visible/selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Duplicate the sheet; Excel auto-names the copy "machine (2)" and
#    inserts it immediately after the source sheet.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# 2. Delete the "ASHE2" row (row 3) on the original sheet; "ATLAS" moves
#    up to become row 3 and the sheet's dimension shrinks to A1:H3.
$ws1.Rows.Item(3).Delete()

# 3. Swap the CalendarioId string for the two remaining machine rows.
$ws1.Range("H2").Value = "CAL-24x5"
$ws1.Range("H3").Value = "CAL-24x5"

# 4. Both tabs end up with the whole of row 3 selected; re-select row 3 on
#    the copy too, then put the selection/active state back on the first
#    sheet (which stays the visible tab).
$ws2.Rows.Item(3).Select()
$ws1.Rows.Item(3).Select()
$ws1.Activate()
